$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rebalance_localizations")

# Insert 3 new rows for the "description" keys (they sort alphabetically right
# before gui/menu/research/description/fire_control_station, which currently
# sits at row 335).
$ws.Rows.Item(335).Insert()
$ws.Rows.Item(335).Insert()
$ws.Rows.Item(335).Insert()

# Insert 3 new rows for the "name" keys (they sort alphabetically right before
# gui/menu/research/name/fire_control_station, which after the inserts above
# now sits at row 506).
$ws.Rows.Item(506).Insert()
$ws.Rows.Item(506).Insert()
$ws.Rows.Item(506).Insert()

# Populate column A (the localization keys) first, in the order the keys were
# authored, so the shared-string table records them in that sequence.
$ws.Range("A506").Value = "gui/menu/research/name/deep_ionization"
$ws.Range("A507").Value = "gui/menu/research/name/deep_ionization_lvl_2"
$ws.Range("A508").Value = "gui/menu/research/name/deep_ionization_lvl_3"
$ws.Range("A335").Value = "gui/menu/research/description/deep_ionization"
$ws.Range("A336").Value = "gui/menu/research/description/deep_ionization_lvl_2"
$ws.Range("A337").Value = "gui/menu/research/description/deep_ionization_lvl_3"

# Now populate column B (the English text), again in authoring order.
$ws.Range("B506").Value = "Deep ionization"
$ws.Range("B507").Value = "Advanced deep ionization"
$ws.Range("B508").Value = "Extreme ionization"
$ws.Range("B335").Value = "High energy physics allows the production of deeply ionized plasmas with strong electrical charge"
$ws.Range("B336").Value = "Advanced techiques allow production of even deeper ionized plasmas carriying extreme charge"
$ws.Range("B337").Value = "Production of superhot plasmas with total atomic ionization. Any more energy and the plasma becomes instable with fusion reactions emerging."

# Apply the column-B style (s="1", used throughout column B) to the new cells.
$ws.Range("B335:B337").Style = $ws.Range("B338").Style
$ws.Range("B506:B508").Style = $ws.Range("B509").Style

# Restore the sort-range / selection bookkeeping to match the post-edit sheet.
$ws.Range("B24").Select()
